$wb = $excel.ActiveWorkbook

# --- Typography sheet ---
$wsTypo = $wb.Worksheets.Item("Typography")

# Update the size of the "Small" typography entry (row 4) from 20 to 30
$wsTypo.Range("D4").Value = 30

# Fill in the new wildcard / ellipsis related columns for rows 4-6
$wsTypo.Range("G4").Value = "."
$wsTypo.Range("H4").Value = ""
$wsTypo.Range("I4").Value = "a-z,A-Z,0-9"
$wsTypo.Range("J4").Value = ""

$wsTypo.Range("G5").Value = ""
$wsTypo.Range("H5").Value = ""
$wsTypo.Range("I5").Value = ""
$wsTypo.Range("J5").Value = ""

$wsTypo.Range("G6").Value = ""
$wsTypo.Range("H6").Value = ""
$wsTypo.Range("I6").Value = ""
$wsTypo.Range("J6").Value = ""

# --- Translation sheet ---
$wsTrans = $wb.Worksheets.Item("Translation")

# Add two new translation rows for the LPS25HB sensor readings (temperature and pressure)
$wsTrans.Range("B4").Value = "SingleUseId1"
$wsTrans.Range("C4").Value = "Default"
$wsTrans.Range("D4").Value = "Left"
$wsTrans.Range("E4").Value = "LTR"
$wsTrans.Range("F4").Value = "Temperatura: <value> *C"

$wsTrans.Range("B5").Value = "SingleUseId2"
$wsTrans.Range("C5").Value = "Default"
$wsTrans.Range("D5").Value = "Left"
$wsTrans.Range("E5").Value = "LTR"
$wsTrans.Range("F5").Value = "Cisnienie: <value> hpa"
